# Scheduled-runner refresh of Leve market-price snapshots (currentAveragePrice*
# and the derived Leve profit columns H:N) across all eight crafting-job sheets.
# Values below are pulled from the upstream market-board pull described in the
# commit; cells that become (or stop being) populated are cleared via "".
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1716
$ws.Range("I18").Value = 621.3333
$ws.Range("K18").Value = 621.3333
$ws.Range("M18").Value = -337.3333
$ws.Range("H40").Value = 6125.3335
$ws.Range("I40").Value = 1485.1666
$ws.Range("J40").Value = 10765.5
$ws.Range("K40").Value = 1485.1666
$ws.Range("L40").Value = 10765.5
$ws.Range("M40").Value = -1310.1666
$ws.Range("N40").Value = -11115.5
$ws.Range("H116").Value = 2866
$ws.Range("J116").Value = 2866
$ws.Range("L116").Value = 2866
$ws.Range("N116").Value = -9750
$ws.Range("H125").Value = 2020.8572
$ws.Range("I125").Value = 650
$ws.Range("J125").Value = 2249.3333
$ws.Range("K125").Value = 5850
$ws.Range("L125").Value = 20243.9997
$ws.Range("M125").Value = -3390
$ws.Range("N125").Value = -25163.9997
$ws.Range("H127").Value = 1003.4
$ws.Range("J127").Value = 2266.6667
$ws.Range("L127").Value = 6800.000100000001
$ws.Range("N127").Value = -16720.0001
$ws.Range("H137").Value = 3049.6775
$ws.Range("I137").Value = 2549.6365
$ws.Range("K137").Value = 7648.9095
$ws.Range("M137").Value = -5098.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4804.146
$ws.Range("I32").Value = 3926.1086
$ws.Range("K32").Value = 3926.1086
$ws.Range("M32").Value = -3639.1086
$ws.Range("H39").Value = 0
$ws.Range("I39").Value = 0
$ws.Range("K39").Value = 0
$ws.Range("M39").Value = ""
$ws.Range("H63").Value = 6906.0713
$ws.Range("I63").Value = 3448.3333
$ws.Range("K63").Value = 3448.3333
$ws.Range("M63").Value = -2762.3333
$ws.Range("H66").Value = 6906.0713
$ws.Range("I66").Value = 3448.3333
$ws.Range("K66").Value = 17241.6665
$ws.Range("M66").Value = -13809.6665
$ws.Range("H97").Value = 1143.7273
$ws.Range("I97").Value = 1143.7273
$ws.Range("K97").Value = 1143.7273
$ws.Range("M97").Value = -647.7273
$ws.Range("H125").Value = 1000000
$ws.Range("J125").Value = 1000000
$ws.Range("L125").Value = 1000000
$ws.Range("N125").Value = -1009840

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("I20").Value = 4534.5264
$ws.Range("J20").Value = 4957
$ws.Range("K20").Value = 4534.5264
$ws.Range("L20").Value = 4957
$ws.Range("M20").Value = -4287.5264
$ws.Range("N20").Value = -5451
$ws.Range("H22").Value = 924.1667
$ws.Range("J22").Value = 1733.3334
$ws.Range("L22").Value = 1733.3334
$ws.Range("N22").Value = -2079.3334
$ws.Range("H94").Value = 2148.375
$ws.Range("I94").Value = 1826.7142
$ws.Range("K94").Value = 1826.7142
$ws.Range("M94").Value = -1375.7142
$ws.Range("H105").Value = 13435.034
$ws.Range("I105").Value = 11884.55
$ws.Range("K105").Value = 11884.55
$ws.Range("M105").Value = -10137.55

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").Value = ""
$ws.Range("H86").Value = 4949
$ws.Range("I86").Value = 4623.5
$ws.Range("J86").Value = 5600
$ws.Range("K86").Value = 4623.5
$ws.Range("L86").Value = 5600
$ws.Range("M86").Value = -3500.5
$ws.Range("N86").Value = -7846
$ws.Range("H89").Value = 4949
$ws.Range("I89").Value = 4623.5
$ws.Range("J89").Value = 5600
$ws.Range("K89").Value = 23117.5
$ws.Range("L89").Value = 28000
$ws.Range("M89").Value = -17501.5
$ws.Range("N89").Value = -39232
$ws.Range("H107").Value = 1434.4762
$ws.Range("I107").Value = 1212.8667
$ws.Range("K107").Value = 1212.8667
$ws.Range("M107").Value = 707.1333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H119").Value = 13216.5
$ws.Range("I119").Value = 700
$ws.Range("J119").Value = 15004.571
$ws.Range("K119").Value = 2100
$ws.Range("L119").Value = 45013.713
$ws.Range("M119").Value = 2738
$ws.Range("N119").Value = -54689.713
$ws.Range("H120").Value = 5000
$ws.Range("I120").Value = 5000
$ws.Range("K120").Value = 15000
$ws.Range("M120").Value = -10162
$ws.Range("H129").Value = 20844962
$ws.Range("J129").Value = 23822528
$ws.Range("L129").Value = 71467584
$ws.Range("N129").Value = -71477584
$ws.Range("H140").Value = 2644.1482
$ws.Range("I140").Value = 1675.7
$ws.Range("J140").Value = 5411.143
$ws.Range("K140").Value = 5027.1
$ws.Range("L140").Value = 16233.429
$ws.Range("M140").Value = 152.8999999999996
$ws.Range("N140").Value = -26593.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3651.1
$ws.Range("I97").Value = 2948.1667
$ws.Range("J97").Value = 4705.5
$ws.Range("K97").Value = 2948.1667
$ws.Range("L97").Value = 4705.5
$ws.Range("M97").Value = -2452.1667
$ws.Range("N97").Value = -5697.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 8448.615
$ws.Range("I46").Value = 5207.75
$ws.Range("J46").Value = 9889
$ws.Range("K46").Value = 5207.75
$ws.Range("L46").Value = 9889
$ws.Range("M46").Value = -5019.75
$ws.Range("N46").Value = -10265
$ws.Range("H116").Value = 184772.83
$ws.Range("I116").Value = 50000
$ws.Range("K116").Value = 50000
$ws.Range("M116").Value = -45411

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 105165.664
$ws.Range("J63").Value = 105165.664
$ws.Range("L63").Value = 105165.664
$ws.Range("N63").Value = -106413.664
$ws.Range("H64").Value = 21648
$ws.Range("I64").Value = 14077
$ws.Range("J64").Value = 31742.666
$ws.Range("K64").Value = 14077
$ws.Range("L64").Value = 31742.666
$ws.Range("M64").Value = -13829
$ws.Range("N64").Value = -32238.666
$ws.Range("H66").Value = 105165.664
$ws.Range("J66").Value = 105165.664
$ws.Range("L66").Value = 315496.992
$ws.Range("N66").Value = -321736.992
$ws.Range("H67").Value = 21648
$ws.Range("I67").Value = 14077
$ws.Range("J67").Value = 31742.666
$ws.Range("K67").Value = 14077
$ws.Range("L67").Value = 31742.666
$ws.Range("M67").Value = -13219
$ws.Range("N67").Value = -33458.666
$ws.Range("H75").Value = 49995
$ws.Range("I75").Value = 49995
$ws.Range("K75").Value = 49995
$ws.Range("M75").Value = -49059
$ws.Range("H78").Value = 49995
$ws.Range("I78").Value = 49995
$ws.Range("K78").Value = 149985
$ws.Range("M78").Value = -145305
$ws.Range("H133").Value = 82357.5
$ws.Range("J133").Value = 82357.5
$ws.Range("L133").Value = 82357.5
$ws.Range("N133").Value = -92477.5
